$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 74
$ws.Range("I2").Value = 142
$ws.Range("J2").Value = 645
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 165
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = 111
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 2
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 71
$ws.Range("T2").Value = 106
$ws.Range("V2").Value = 983
$ws.Range("X2").Value = 1036
$ws.Range("Z2").Value = 10
$ws.Range("AA2").Value = 4
